$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cells I1 and J1 (new columns "I0" and "IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the formatting (bold, centered, bordered) from the existing H1 header
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Data values for columns I (I0) and J (IF), rows 2-26
$data = @(
    @(5, 7),
    @(8, 8),
    @(7, 8),
    @(4, 5),
    @(9, 9),
    @(6, 6),
    @(7, 7),
    @(8, 8),
    @(6, 8),
    @(8, 8),
    @(5, 6),
    @(8, 8),
    @(1, 3),
    @(7, 7),
    @(9, 9),
    @(7, 8),
    @(8, 8),
    @(4, 5),
    @(6, 6),
    @(3, 6),
    @(5, 5),
    @(3, 4),
    @(5, 5),
    @(5, 6),
    @(3, 3)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
    $row++
}
